# Kapeleshh CV Purple — "Change CG and Update Preprint"
#
# 1. CGPA table cell: 7.31/10 -> 7.78/10
# 2. Add a new bulleted paragraph under the CKD / CNN project bullet:
#      "Published preprint of the paper on public health approach to
#       prediction of CKD"
#    with "preprint" hyperlinked out to the published preprint.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Update the CGPA value in the education table.
# ---------------------------------------------------------------------
$tableShape = $s.Shapes.Item(3)
$tbl = $tableShape.Table
$cgpaCell = $tbl.Cell(2, 3)
$cgpaRange = $cgpaCell.Shape.TextFrame.TextRange
$cgpaRange.Text = $cgpaRange.Text.Replace("7.31/10", "7.78/10")

# ---------------------------------------------------------------------
# 2) Append the "Published preprint ..." bullet to the CKD project shape.
# ---------------------------------------------------------------------
$projShape = $s.Shapes.Item(7)
$tr = $projShape.TextFrame.TextRange

$startLen = $tr.Text.Length

$part1 = "Published "
$part2 = "preprint"
$part3 = " of the paper on public health approach to prediction of CKD"

# Inserting a leading carriage return starts a brand-new paragraph that
# inherits the bullet formatting of the preceding paragraph.
$tr.InsertAfter("`r" + $part1 + $part2 + $part3) | Out-Null

$part2Start = $startLen + 1 + $part1.Length + 1
$preprintRange = $tr.Characters($part2Start, $part2.Length)

$actionSetting = $preprintRange.ActionSettings(1)
$actionSetting.Hyperlink.Address = "https://www.medrxiv.org/content/10.1101/2020.06.18.20134304v1"
